# Insert a new weekly price record for "Camote" (Vega Modelo de Temuco)
# into the data table. The new record is inserted as row 61, pushing all
# subsequent rows down by one (dimension grows from A1:R147 to A1:R148).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 61
$ws.Rows("61:61").Insert()

$ws.Cells.Item($newRow, 1).Value  = 10
$ws.Cells.Item($newRow, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item($newRow, 3).Value  = "La Araucanía"
$ws.Cells.Item($newRow, 4).Value  = 44895
$ws.Cells.Item($newRow, 5).Value  = 9
$ws.Cells.Item($newRow, 6).Value  = 100114002
$ws.Cells.Item($newRow, 7).Value  = "Camote"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 15
$ws.Cells.Item($newRow, 11).Value = 22000
$ws.Cells.Item($newRow, 12).Value = 22000
$ws.Cells.Item($newRow, 13).Value = 22000
$ws.Cells.Item($newRow, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item($newRow, 15).Value = "Perú"
$ws.Cells.Item($newRow, 16).Value = 1100
$ws.Cells.Item($newRow, 17).Value = 20
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
